# Add a new row of data (row 12) to the worksheet, matching the style/format
# used by the existing data rows, with all three cells stored as plain text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target cells to be treated as Text so that the phone number,
# DDD code and date are stored literally (as strings) instead of being
# auto-converted to a number / date by Excel's input parsing.
$ws.Range("A12:C12").NumberFormat = "@"

$ws.Range("A12").Value = "+556298529715"
$ws.Range("B12").Value = "62"
$ws.Range("C12").Value = "2024-07-09"

# Copy the formatting (style) of the previous data row (row 11) onto the
# new row so the new cells share the same cell style index as the other
# data rows, then restore the General number format on top of it.
$ws.Range("A11:C11").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)
